$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for column C: "ano" (year)
$ws.Range("C1").Value = "ano"

# Fill C2:C36 with the year value 2023
$ws.Range("C2:C36").Value = 2023
